$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "7. Gestione del File System",
    "8. Gestione della Memoria e memoria virtuale",
    "9. Processi e Risorse, Calcolo Parallelo e Distribuito",
    "10. La concorrenza: Il modello a processi; Processi sequenziali e paralleli, I Thread",
    "11. Elaborazione sequenziale e concorrente, La Concorrenza tra processi",
    "12. Scheduling della CPU e Algoritmi di scheduling: sincronizzazione e comunicazione tra processi",
    "13. Sezione critica e mutua esclusione",
    "14. I semafori, produttori/consumatori, lettori/scrittori, deadlock, monitor e scambio di messaggi",
    "15. Comandi fondamentali di Linux e bash scripting",
    "16. Comandi fondamentali di dos/PowerShell scripting"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 9 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

$ws.Range("H24").Select()
